# Address label and entering certificate details
# Insert a new "RobotExceptionFolder" row into the ROBOT PARAMETERS table on
# Sheet1, right after "RobotWorkingFolder" (row 25) and before
# "EmailApplicationsFolder" (previously row 26). This pushes every
# subsequent row down by one and grows Table1 by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a blank row at row 26 - everything from the old row 26 downward
# shifts down to row 27 onward, carrying its formatting (styles / row
# heights) with it.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row with the new parameter.
$ws.Cells.Item(26, 1).Value = "RobotExceptionFolder"
$ws.Cells.Item(26, 2).Value = "\Desktop\RobotWorkingFolder\Exceptions\"
$ws.Cells.Item(26, 3).Value = "Folder path for exceptions folder"

# Grow the Table1 ListObject so its range / autofilter cover the new row.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:C37"))

# Match the author's final cursor position/selection in the saved file.
$ws.Activate()
$ws.Range("C27").Select()
